$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 52 (Giorno_Cont 50) ---
$ws.Cells.Item(52, 3).Value = 15113
$ws.Cells.Item(52, 4).Value = 1016

# --- New rows 120-133 (Giorno_Cont 118-131, dates 5/19/20 - 6/1/20) ---
$newRows = @(
    @(118, "5/19/20", 226699, 32169, 129401),
    @(119, "5/20/20", 227364, 32330, 132282),
    @(120, "5/21/20", 228006, 32486, 134560),
    @(121, "5/22/20", 228658, 32616, 136720),
    @(122, "5/23/20", 229327, 32735, 138840),
    @(123, "5/24/20", 229858, 32785, 140479),
    @(124, "5/25/20", 230158, 32877, 141981),
    @(125, "5/26/20", 230555, 32955, 144658),
    @(126, "5/27/20", 231139, 33072, 147101),
    @(127, "5/28/20", 231732, 33142, 150604),
    @(128, "5/29/20", 232248, 33229, 152844),
    @(129, "5/30/20", 232664, 33340, 155633),
    @(130, "5/31/20", 232997, 33415, 157507),
    @(131, "6/1/20", 233197, 33475, 158355)
)

# Helper cell used to stage a text formula result for the date strings so
# that the B-column cells receive a genuine text value (matching the
# existing rows) instead of Excel's automatic date-literal conversion,
# and without registering any new number-format styles.
$helper = $ws.Cells.Item(1000, 100)

$r = 120
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]

    $helper.Formula = '="' + $row[1] + '"'
    $helper.Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4163)

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$helper.ClearContents()
